$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Range("H17").Value = 527089.5
$ws.Range("J17").Value = 577253.4399999999
$ws.Range("L17").Value = 1731760.32
$ws.Range("N17").Value = -1732096.32

# Row 132 (ALC)
$ws.Range("H132").Value = 8189.2
$ws.Range("I132").Value = 8189.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24567.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -22037.6
$ws.Range("N132").ClearContents()

# Row 135 (ALC)
$ws.Range("H135").Value = 115385720
$ws.Range("I135").Value = 100000690
$ws.Range("K135").Value = 900006210
$ws.Range("M135").Value = -900003675

# Row 137 (ALC)
$ws.Range("H137").Value = 2654.3
$ws.Range("I137").Value = 2304.1428
$ws.Range("J137").Value = 3471.3333
$ws.Range("K137").Value = 6912.428400000001
$ws.Range("L137").Value = 10413.9999
$ws.Range("M137").Value = -4362.428400000001
$ws.Range("N137").Value = -15513.9999

# Row 138 (ALC)
$ws.Range("H138").Value = 5682.304
$ws.Range("I138").Value = 3917.1
$ws.Range("J138").Value = 7040.154
$ws.Range("K138").Value = 11751.3
$ws.Range("L138").Value = 21120.462
$ws.Range("M138").Value = -6611.299999999999
$ws.Range("N138").Value = -31400.462

# Row 141 (ALC)
$ws.Range("H141").Value = 2416.6
$ws.Range("I141").Value = 2716.25
$ws.Range("K141").Value = 8148.75
$ws.Range("M141").Value = -2968.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3147.261
$ws.Range("I32").Value = 1506.3954
$ws.Range("K32").Value = 1506.3954
$ws.Range("M32").Value = -1219.3954

# Row 36 (ARM)
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

# Row 61 (ARM)
$ws.Range("H61").Value = 90910950
$ws.Range("I61").Value = 90910950
$ws.Range("K61").Value = 90910950
$ws.Range("M61").Value = -90910738

# Row 74 (ARM)
$ws.Range("H74").Value = 29416166
$ws.Range("I74").Value = 30307262
$ws.Range("K74").Value = 30307262
$ws.Range("M74").Value = -30306388

# Row 77 (ARM)
$ws.Range("H77").Value = 29416166
$ws.Range("I77").Value = 30307262
$ws.Range("K77").Value = 151536310
$ws.Range("M77").Value = -151531942

# Row 132 (ARM)
$ws.Range("H132").Value = 6253288.5
$ws.Range("I132").Value = 9094339
$ws.Range("K132").Value = 27283017
$ws.Range("M132").Value = -27280487

# Row 136 (ARM)
$ws.Range("H136").Value = 90910950
$ws.Range("I136").Value = 90910950
$ws.Range("K136").Value = 272732850
$ws.Range("M136").Value = -272730300

$ws = $wb.Worksheets.Item("BSM")
# Row 102 (BSM)
$ws.Range("H102").Value = 8217
$ws.Range("I102").Value = 8217
$ws.Range("K102").Value = 8217
$ws.Range("M102").Value = -4972

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 602.05
$ws.Range("J22").Value = 499.5
$ws.Range("L22").Value = 499.5
$ws.Range("N22").Value = -1199.5

# Row 31 (CRP)
$ws.Range("H31").Value = 3807.9773
$ws.Range("I31").Value = 6975.5386
$ws.Range("J31").Value = 2479.6453
$ws.Range("K31").Value = 6975.5386
$ws.Range("L31").Value = 2479.6453
$ws.Range("M31").Value = -6680.5386
$ws.Range("N31").Value = -3069.6453

# Row 34 (CRP)
$ws.Range("H34").Value = 3807.9773
$ws.Range("I34").Value = 6975.5386
$ws.Range("J34").Value = 2479.6453
$ws.Range("K34").Value = 6975.5386
$ws.Range("L34").Value = 2479.6453
$ws.Range("M34").Value = -6773.5386
$ws.Range("N34").Value = -2883.6453

# Row 58 (CRP)
$ws.Range("H58").Value = 25006354
$ws.Range("I58").Value = 26322214
$ws.Range("J58").Value = 5014
$ws.Range("K58").Value = 26322214
$ws.Range("L58").Value = 5014
$ws.Range("M58").Value = -26322011
$ws.Range("N58").Value = -5420

# Row 132 (CRP)
$ws.Range("H132").Value = 31255002
$ws.Range("I132").Value = 38466544
$ws.Range("J132").Value = 4992.1665
$ws.Range("K132").Value = 115399632
$ws.Range("L132").Value = 14976.4995
$ws.Range("M132").Value = -115397102
$ws.Range("N132").Value = -20036.4995

# Row 134 (CRP)
$ws.Range("H134").Value = 11908904
$ws.Range("I134").Value = 12503849
$ws.Range("K134").Value = 37511547
$ws.Range("M134").Value = -37509012

# Row 136 (CRP)
$ws.Range("H136").Value = 25006354
$ws.Range("I136").Value = 26322214
$ws.Range("J136").Value = 5014
$ws.Range("K136").Value = 78966642
$ws.Range("L136").Value = 15042
$ws.Range("M136").Value = -78964092
$ws.Range("N136").Value = -20142

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (CUL)
$ws.Range("H68").Value = 6720.2646
$ws.Range("I68").Value = 1666
$ws.Range("J68").Value = 6953.5386
$ws.Range("K68").Value = 4998
$ws.Range("L68").Value = 20860.6158
$ws.Range("M68").Value = -4187
$ws.Range("N68").Value = -22482.6158

# Row 71 (CUL)
$ws.Range("H71").Value = 6720.2646
$ws.Range("I71").Value = 1666
$ws.Range("J71").Value = 6953.5386
$ws.Range("K71").Value = 14994
$ws.Range("L71").Value = 62581.8474
$ws.Range("M71").Value = -10938
$ws.Range("N71").Value = -70693.8474

# Row 107 (CUL)
$ws.Range("H107").Value = 2033.3334
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 900
$ws.Range("M107").Value = 1020

# Row 129 (CUL)
$ws.Range("H129").Value = 1700.92
$ws.Range("J129").Value = 3905.3333
$ws.Range("L129").Value = 11715.9999
$ws.Range("N129").Value = -21715.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 74 (GSM)
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77 (GSM)
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 107 (GSM)
$ws.Range("H107").Value = 2086.5715
$ws.Range("I107").Value = 1575.75
$ws.Range("J107").Value = 5151.5
$ws.Range("K107").Value = 1575.75
$ws.Range("L107").Value = 5151.5
$ws.Range("M107").Value = 344.25
$ws.Range("N107").Value = -8991.5

# Row 132 (GSM)
$ws.Range("H132").Value = 12503809
$ws.Range("I132").Value = 17861188
$ws.Range("J132").Value = 3261.3333
$ws.Range("K132").Value = 53583564
$ws.Range("L132").Value = 9783.999899999999
$ws.Range("M132").Value = -53581034
$ws.Range("N132").Value = -14843.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (LTW)
$ws.Range("H68").Value = 1324151.2
$ws.Range("I68").Value = 2636103.5
$ws.Range("J68").Value = 12199.2
$ws.Range("K68").Value = 2636103.5
$ws.Range("L68").Value = 12199.2
$ws.Range("M68").Value = -2635354.5
$ws.Range("N68").Value = -13697.2

# Row 71 (LTW)
$ws.Range("H71").Value = 1324151.2
$ws.Range("I71").Value = 2636103.5
$ws.Range("J71").Value = 12199.2
$ws.Range("K71").Value = 13180517.5
$ws.Range("L71").Value = 60996
$ws.Range("M71").Value = -13176773.5
$ws.Range("N71").Value = -68484

# Row 132 (LTW)
$ws.Range("H132").Value = 10003813
$ws.Range("I132").Value = 11367372
$ws.Range("J132").Value = 4381.3335
$ws.Range("K132").Value = 34102116
$ws.Range("L132").Value = 13144.0005
$ws.Range("M132").Value = -34099586
$ws.Range("N132").Value = -18204.0005

# Row 141 (LTW)
$ws.Range("H141").Value = 99995
$ws.Range("J141").Value = 99995
$ws.Range("L141").Value = 99995
$ws.Range("N141").Value = -110355

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Range("H96").Value = 1499.6666
$ws.Range("I96").Value = 1499.6666
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1499.6666
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -126.6666
$ws.Range("N96").ClearContents()

# Row 132 (WVR)
$ws.Range("H132").Value = 45474330
$ws.Range("I132").Value = 62507948
$ws.Range("K132").Value = 187523844
$ws.Range("M132").Value = -187521314
